# Update countries & provincias Spain
# Applies the data refresh captured in the commit's diff of paises.xlsx:
#   - swap the "Santa Lucia" / "Nueva Caledonia" rows (shared-string reorder
#     in the XML manifests as these two adjacent rows trading display text)
#   - bump the "Datos actualizados..." timestamp in A1
#   - refresh the numeric counters for a handful of country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- timestamp in A1 -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 06:13"

# --- Santa Lucia / Nueva Caledonia swap (rows 207 & 208) -------------
$ws.Cells.Item(207, 1).Value = "Santa Lucia"
$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"

# --- numeric refresh for countries with updated counts ---------------
# row -> @{ col = value }
$updates = @{
    5   = @{ 2 = 6312584; 3 = 2317; 4 = 5273201; 5 = 940675 }
    23  = @{ 2 = 312806;  3 = 543;  4 = 297497;  5 = 8825;  7 = 5;  8 = 6484 }
    36  = @{ 2 = 118452;  3 = 1337; 4 = 19454;   5 = 88982; 7 = 15; 8 = 10016 }
    39  = @{ 2 = 107979;  3 = 71;   4 = 102937;  5 = 3317 }
    135 = @{ 2 = 4123;    3 = 101;  4 = 2206;    5 = 1821; 7 = 1;  8 = 96 }
    141 = @{ 2 = 3569;    3 = 5;    4 = 3379 }
    157 = @{ 2 = 1992;    3 = 49;   4 = 1246;    5 = 719;  7 = 1;  8 = 27 }
    186 = @{ 4 = 307;     5 = 6 }
    187 = @{ 2 = 282;     3 = 1;    4 = 222;     5 = 60 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item([int]$row, [int]$col).Value = $cols[$col]
    }
}
